$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert the new "Projects" worksheet right after "References" (so it
#    lands before "DLLs", matching the target sheet order).
# ---------------------------------------------------------------------------
$refSheet = $wb.Worksheets.Item("References")
$projects = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $refSheet)
$projects.Name = "Projects"

$tableDoc = $wb.Worksheets.Item("Table Doc")
$seasonInfo = $wb.Worksheets.Item("SeasonInfo")

# ---------------------------------------------------------------------------
# 2. Update the "Table Doc" header cell C1: new text + new style (bordered,
#    bold, yellow fill, vertical-top/wrap-text alignment). We build the style
#    by copying an existing bold/yellow/bordered cell and then layering the
#    alignment on top so the engine reuses/creates the exact xf combination.
# ---------------------------------------------------------------------------
$styleSrc = $seasonInfo.Range("A1")
$c1 = $tableDoc.Range("C1")
$styleSrc.Copy()
$c1.PasteSpecial(-4122)
$c1.WrapText = $true
$c1.VerticalAlignment = -4160
$c1.Value = "For TMs Analysis " + [char]10 + "Tesing"

# ---------------------------------------------------------------------------
# 3. Populate the "Projects" sheet. Values are written in the exact order
#    the original author typed them so newly-created shared strings line up
#    with the target workbook (existing strings such as "Bball.VbClasses"
#    are simply reused).
# ---------------------------------------------------------------------------
$projects.Range("A1").Value = "Project"
$projects.Range("A3").Value = "Bball.DAL"
$projects.Range("A4").Value = "Bball.DataBaseFunctions"
$projects.Range("A5").Value = "Bball.lBAL"
$projects.Range("A6").Value = "Bball.lDAL"
$projects.Range("A7").Value = "Bball.Unity (unloaded)"
$projects.Range("A8").Value = "Bball.VbClasses"
$projects.Range("A9").Value = "Bball.VbClasseslnterfaces"
$projects.Range("A10").Value = "BballConsole (unloaded)"
$projects.Range("A11").Value = "BballMVC"
$projects.Range("A12").Value = "BballMVC.DTOs"
$projects.Range("A13").Value = "BballMVC.IDTOs"
$projects.Range("A14").Value = "BballMVC.Tests (unloaded)"
$projects.Range("A15").Value = "HtmlParser"
$projects.Range("A16").Value = "HtmlParsing"
$projects.Range("A17").Value = "StringExtensionMethods"
$projects.Range("A18").Value = "SysDAL.Functions"
$projects.Range("A19").Value = "Trace"
$projects.Range("A20").Value = "TTl.Logger"
$projects.Range("A22").Value = "UnitTest_CS"
$projects.Range("A23").Value = "UnitTestProject1 (unloaded)"
$projects.Range("B1").Value = ".net ver"
$projects.Range("A21").Value = "TTI.Models"
$projects.Range("A2").Value = "Bball.BAL"

# Numeric / text values for column B.
$projects.Range("B2").Value = 472
$projects.Range("B3").Value = 472
$projects.Range("B4").Value = 45
$projects.Range("B5").Value = 472
$projects.Range("B6").Value = 45
$projects.Range("B7").Value = "x"
$projects.Range("B8").Value = 45
$projects.Range("B9").Value = 45
$projects.Range("B10").Value = "x"
$projects.Range("B11").Value = 472
$projects.Range("B12").Value = 45
$projects.Range("B13").Value = 45
$projects.Range("B14").Value = "x"
$projects.Range("B15").Value = 45
$projects.Range("B16").Value = 45
$projects.Range("B17").Value = 45
$projects.Range("B18").Value = 45
$projects.Range("B19").Value = 45
$projects.Range("B20").Value = 45
$projects.Range("B21").Value = 45
$projects.Range("B22").Value = 472

# Header style for A1: same bold/yellow/filled font+fill as other section
# headers but without the border (layered on top of a matching base cell so
# the resulting xf exactly matches font+fill only).
$headerStyleSrc = $seasonInfo.Range("A1")
$headerStyleSrc.Copy()
$projects.Range("A1").PasteSpecial(-4122)
$projects.Range("A1").Borders.LineStyle = -4142

# Column widths (approximate AutoFit result for the two columns).
$projects.Columns.Item(1).ColumnWidth = 27.57
$projects.Columns.Item(2).ColumnWidth = 7.86

# Sheet view: freeze nothing, just scroll/select like the source file.
$projects.Range("A3").Select()

# Page setup to match target (portrait orientation).
$projects.PageSetup.Orientation = 1

# ---------------------------------------------------------------------------
# 4. "Table Doc": clear the old bottom-right selection anchor and move it to
#    C2 (matches the diff's updated <selection pane="bottomRight" .../>).
# ---------------------------------------------------------------------------
$tableDoc.Activate()
$tableDoc.Range("C2").Select()

# ---------------------------------------------------------------------------
# 5. "References": scroll the frozen pane back up to the top (A2) instead of
#    A26.
# ---------------------------------------------------------------------------
$references = $wb.Worksheets.Item("References")
$references.Activate()
$winRef = $excel.ActiveWindow
$winRef.ScrollRow = 2

# ---------------------------------------------------------------------------
# 6. Finally make "Projects" the active sheet/tab (matches activeTab update).
# ---------------------------------------------------------------------------
$projects.Activate()
$projects.Range("A3").Select()

Write-Host "Edit complete"
